$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New task names first, in row order, so the shared-string table picks
#     up "impl message properties", "template.jsp", "register",
#     "Theater/play simple" and "login" in that order before any new date
#     ("18.04.2016") or status ("Inprocess") string is introduced.
$ws.Range("B9").Value = "impl message properties"
$ws.Range("B10").Value = "template.jsp"
$ws.Range("B11").Value = "register"
$ws.Range("B12").Value = "Theater/play simple"
$ws.Range("B13").Value = "login"

# --- Row 6: E6 keeps the existing "16.04.2016" text (its shared-string
# index will simply shift once "Pending" becomes unused below).
$ws.Range("E6").Value = "16.04.2016"

# --- Row 7: D7 and E7 stay "16.04.2016"
$ws.Range("D7").Value = "16.04.2016"
$ws.Range("E7").Value = "16.04.2016"

# --- Row 8: impl hibernate model task now has hours=1, dates=16.04.2016, status=Done
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = "16.04.2016"
$ws.Range("E8").Value = "16.04.2016"
$ws.Range("F8").Value = "Done"

# --- Row 9: impl message properties
$ws.Range("C9").Value = 0.5
$ws.Range("D9").Value = "16.04.2016"
$ws.Range("E9").Value = "16.04.2016"
$ws.Range("F9").Value = "Done"

# --- Row 10: template.jsp
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = "16.04.2016"
$ws.Range("E10").Value = "16.04.2016"
$ws.Range("F10").Value = "Done"

# --- Row 11: register (introduces the new "18.04.2016" string)
$ws.Range("C11").Value = 6
$ws.Range("D11").Value = "16.04.2016"
$ws.Range("E11").Value = "18.04.2016"
$ws.Range("F11").Value = "Done"

# --- Row 12: Theater/play simple
$ws.Range("C12").Value = 0.2
$ws.Range("D12").Value = "16.04.2016"
$ws.Range("E12").Value = "16.04.2016"
$ws.Range("F12").Value = "Done"

# --- Row 13: login (introduces the new "Inprocess" string last)
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = "18.04.2016"
$ws.Range("F13").Value = "Inprocess"

# --- Data validation sqref had a redundant "F1" entry which should be
#     removed; recreate the validation on F1:F1048576 only.
$ws.Range("F1:F1048576").Validation.Delete()
$ws.Range("F1:F1048576").Validation.Add(3, 1, 1, '"Pending,Inprocess, Done"')
$ws.Range("F1:F1048576").Validation.InCellDropdown = $true
$ws.Range("F1:F1048576").Validation.IgnoreBlank = $true

# --- Update the current selection to match the target view state
$ws.Range("F15").Select()

$wb.Save()
